$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Testing plans score was left blank before; set it to the maximum (3),
# same as this section's other fully-scored rows.
$ws.Range("D44").Value = 3

# Move the view / selection to where the edit was made, mirroring the
# author's final cursor position in the saved workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C43:D43").Select()
